$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.23018405348877
$ws.Range("C2").Value = 11.44452154835927
$ws.Range("E2").Value = 12.33771884286134
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 22.95969453361531
$ws.Range("H2").Value = 12.80026029692002
$ws.Range("L2").Value = 9.493625208932086
$ws.Range("N2").Value = 16.7289679979256
$ws.Range("O2").Value = 18.71072748547778
$ws.Range("B3").Value = 13.70659715882156
$ws.Range("C3").Value = 11.3741243488934
$ws.Range("E3").Value = 12.37880254957037
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 22.97573783868756
$ws.Range("H3").Value = 12.84617873965466
$ws.Range("L3").Value = 9.465950482281782
$ws.Range("N3").Value = 16.76114993473243
$ws.Range("O3").Value = 18.77689615497752
$ws.Range("B4").Value = 13.37595319360735
$ws.Range("C4").Value = 11.33102180040358
$ws.Range("E4").Value = 12.40640735468651
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 22.99532817282422
$ws.Range("H4").Value = 12.87672416326989
$ws.Range("L4").Value = 9.450517175909466
$ws.Range("N4").Value = 16.78271443504787
$ws.Range("O4").Value = 18.8223937321975
$ws.Range("B5").Value = 13.2391084544898
$ws.Range("C5").Value = 11.3134977869644
$ws.Range("E5").Value = 12.41825424309139
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 23.00575138542714
$ws.Range("H5").Value = 12.8897624425599
$ws.Range("L5").Value = 9.444624618453009
$ws.Range("N5").Value = 16.79195661868225
$ws.Range("O5").Value = 18.84215496447876
$ws.Range("B6").Value = 13.21626458245231
$ws.Range("C6").Value = 11.31059068278503
$ws.Range("E6").Value = 12.42025749265411
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 23.00762920941056
$ws.Range("H6").Value = 12.89196310578516
$ws.Range("L6").Value = 9.443670246993795
$ws.Range("N6").Value = 16.79351874529294
$ws.Range("O6").Value = 18.84550990513092
$ws.Range("B7").Value = 13.37411590168558
$ws.Range("C7").Value = 11.33078528775911
$ws.Range("E7").Value = 12.40656470664603
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 22.99545887827092
$ws.Range("H7").Value = 12.87689761070284
$ws.Range("L7").Value = 9.450436095035759
$ws.Range("N7").Value = 16.78283723732118
$ws.Range("O7").Value = 18.82265530225166
$ws.Range("B8").Value = 14.05166517306092
$ws.Range("C8").Value = 11.42022940474122
$ws.Range("E8").Value = 12.35139035588137
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 22.96320001199311
$ws.Range("H8").Value = 12.81560467114632
$ws.Range("L8").Value = 9.48376223760908
$ws.Range("N8").Value = 16.73969017914451
$ws.Range("O8").Value = 18.732529406857
$ws.Range("B9").Value = 15.29995892473292
$ws.Range("C9").Value = 11.59615082644726
$ws.Range("E9").Value = 12.26209671238581
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 22.97751611594684
$ws.Range("H9").Value = 12.71408753534676
$ws.Range("L9").Value = 9.561268292701286
$ws.Range("N9").Value = 16.66936854077085
$ws.Range("O9").Value = 18.5946016064844
$ws.Range("B10").Value = 16.15911893700896
$ws.Range("C10").Value = 11.7251231913968
$ws.Range("E10").Value = 12.20804927895291
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 23.03556340099611
$ws.Range("H10").Value = 12.65091867708315
$ws.Range("L10").Value = 9.62530891381712
$ws.Range("N10").Value = 16.62637510888124
$ws.Range("O10").Value = 18.51714680488743
$ws.Range("B11").Value = 16.53585040884596
$ws.Range("C11").Value = 11.78360861950835
$ws.Range("E11").Value = 12.18597761336417
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 23.07228748362083
$ws.Range("H11").Value = 12.62466656176653
$ws.Range("L11").Value = 9.655910308661838
$ws.Range("N11").Value = 16.60869097372062
$ws.Range("O11").Value = 18.48713970817006
$ws.Range("B12").Value = 16.67637378170925
$ws.Range("C12").Value = 11.80571817126161
$ws.Range("E12").Value = 12.1779817764356
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 23.08767307300378
$ws.Range("H12").Value = 12.61508326521168
$ws.Range("L12").Value = 9.667702798624337
$ws.Range("N12").Value = 16.60226322992875
$ws.Range("O12").Value = 18.47653175481607
$ws.Range("B13").Value = 16.64620616285577
$ws.Range("C13").Value = 11.80095833177606
$ws.Range("E13").Value = 12.17968770380683
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 23.08429383748368
$ws.Range("H13").Value = 12.61713127557715
$ws.Range("L13").Value = 9.665154095728031
$ws.Range("N13").Value = 16.60363561168001
$ws.Range("O13").Value = 18.47878273361105
$ws.Range("B14").Value = 16.54745474349892
$ws.Range("C14").Value = 11.78542840249905
$ws.Range("E14").Value = 12.18531252823108
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 23.0735236612205
$ws.Range("H14").Value = 12.6238709634327
$ws.Range("L14").Value = 9.656876427459581
$ws.Range("N14").Value = 16.60815677454103
$ws.Range("O14").Value = 18.48625183397657
$ws.Range("B15").Value = 16.48668541807461
$ws.Range("C15").Value = 11.77591064741013
$ws.Range("E15").Value = 12.18880508504825
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 23.06711903664738
$ws.Range("H15").Value = 12.6280458335201
$ws.Range("L15").Value = 9.651832528642521
$ws.Range("N15").Value = 16.61096111243595
$ws.Range("O15").Value = 18.49092530389103
$ws.Range("B16").Value = 16.13420522794535
$ws.Range("C16").Value = 11.72129636849154
$ws.Range("E16").Value = 12.20954237028738
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 23.03337067744108
$ws.Range("H16").Value = 12.65268434940449
$ws.Range("L16").Value = 9.623338023025623
$ws.Range("N16").Value = 16.62756846340544
$ws.Range("O16").Value = 18.51921332361127
$ws.Range("B17").Value = 15.9142757093399
$ws.Range("C17").Value = 11.6877370889052
$ws.Range("E17").Value = 12.22290855854975
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 23.01530733648405
$ws.Range("H17").Value = 12.66843588827845
$ws.Range("L17").Value = 9.606229041547287
$ws.Range("N17").Value = 16.63823603597143
$ws.Range("O17").Value = 18.53790846824122
$ws.Range("B18").Value = 15.78645597438281
$ws.Range("C18").Value = 11.66841794852676
$ws.Range("E18").Value = 12.23083310611219
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 23.00588942322403
$ws.Range("H18").Value = 12.67772949718548
$ws.Range("L18").Value = 9.596527104598602
$ws.Range("N18").Value = 16.64454814743445
$ws.Range("O18").Value = 18.54915326871943
$ws.Range("B19").Value = 15.74295487945605
$ws.Range("C19").Value = 11.66187428325322
$ws.Range("E19").Value = 12.23355685107219
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 23.00286766078116
$ws.Range("H19").Value = 12.6809162785353
$ws.Range("L19").Value = 9.593266220225614
$ws.Range("N19").Value = 16.64671563589548
$ws.Range("O19").Value = 18.55304493353448
$ws.Range("B20").Value = 15.93782524469164
$ws.Range("C20").Value = 11.69131134317514
$ws.Range("E20").Value = 12.22146120611019
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 23.01712967443562
$ws.Range("H20").Value = 12.66673491405017
$ws.Range("L20").Value = 9.608036015841643
$ws.Range("N20").Value = 16.63708220161378
$ws.Range("O20").Value = 18.5358674085587
$ws.Range("B21").Value = 16.57651921084257
$ws.Range("C21").Value = 11.78999102769871
$ws.Range("E21").Value = 12.18365054569341
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 23.07664703528949
$ws.Range("H21").Value = 12.6218816377096
$ws.Range("L21").Value = 9.659302287633091
$ws.Range("N21").Value = 16.60682150753188
$ws.Range("O21").Value = 18.48403745967389
$ws.Range("B22").Value = 16.98144803546714
$ws.Range("C22").Value = 11.8542595990808
$ws.Range("E22").Value = 12.16105066326326
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 23.12416157565933
$ws.Range("H22").Value = 12.59465331516454
$ws.Range("L22").Value = 9.693995828731836
$ws.Range("N22").Value = 16.58861120643608
$ws.Range("O22").Value = 18.45456618946659
$ws.Range("B23").Value = 16.76650535519022
$ws.Range("C23").Value = 11.81998241211365
$ws.Range("E23").Value = 12.17291924811265
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 23.09801602520545
$ws.Range("H23").Value = 12.60899451066858
$ws.Range("L23").Value = 9.6753728542504
$ws.Range("N23").Value = 16.59818722046466
$ws.Range("O23").Value = 18.46989166926246
$ws.Range("B24").Value = 15.92718278765736
$ws.Range("C24").Value = 11.68969550129793
$ws.Range("E24").Value = 12.22211480616378
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 23.01630278317592
$ws.Range("H24").Value = 12.66750318368727
$ws.Range("L24").Value = 9.607218664095422
$ws.Range("N24").Value = 16.637603291939
$ws.Range("O24").Value = 18.53678862466703
$ws.Range("B25").Value = 14.97186388986529
$ws.Range("C25").Value = 11.54856479861774
$ws.Range("E25").Value = 12.284225804324
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 22.96530361275688
$ws.Range("H25").Value = 12.73954796652408
$ws.Range("L25").Value = 9.539030599380323
$ws.Range("N25").Value = 16.68686668350165
$ws.Range("O25").Value = 18.62773631703253
